$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the "Biospecimen identifiers" table, which occupies rows 12-21
# (row 12 is the table header/title row, rows 13-21 are its field rows).
$ws.Rows("12:21").Delete()

# Row deletion does not automatically drop hyperlinks that pointed into the
# removed rows, so clean up any leftover Hyperlink objects whose anchor row
# no longer belongs to the remaining table (rows 1-11).
$changed = $true
while ($changed) {
    $changed = $false
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -gt 11) {
            $hl.Delete()
            $changed = $true
            break
        }
    }
}
